$d = $word.ActiveDocument

$replacements = @(
    @("117÷2=", "487÷7="),
    @("901÷8=", "488÷2="),
    @("184÷3=", "358÷9="),
    @("934÷4=", "495÷6="),
    @("479÷7=", "414÷8="),
    @("493÷7=", "615÷9="),
    @("963÷6=", "749÷4="),
    @("110÷3=", "390÷9="),
    @("463÷2=", "822÷9="),
    @("754÷5=", "786÷8="),
    @("770÷9=", "393÷6="),
    @("669÷2=", "253÷2="),
    @("572÷8=", "774÷3="),
    @("991÷3=", "285÷6="),
    @("687÷2=", "426÷6="),
    @("802÷6=", "360÷6="),
    @("489÷2=", "623÷6="),
    @("277÷9=", "284÷3="),
    @("711÷2=", "112÷3="),
    @("634÷9=", "120÷7="),
    @("338÷3=", "399÷5="),
    @("938÷4=", "613÷5="),
    @("536÷8=", "561÷3="),
    @("346÷4=", "800÷7="),
    @("234÷2=", "862÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
